# Update symbol list values (cryptos.xlsx) per the automated GitHub Actions run.
# Column D holds price-like figures that are stored as TEXT (not numbers) in the
# workbook, so every D-column write below is prefixed with a leading apostrophe
# to force Excel to keep it as literal text (quote-prefixed string) instead of
# silently converting it to a numeric value and losing formatting / precision.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value  = "'245.31"
$ws.Range("D5").Value  = "'0.05768"
$ws.Range("D6").Value  = "'6.458"
$ws.Range("D7").Value  = "'3.149"
$ws.Range("D8").Value  = "'0.8158"
$ws.Range("D9").Value  = "'0.8469"
$ws.Range("D10").Value = "'0.1359"
$ws.Range("D11").Value = "'0.06962"
$ws.Range("D12").Value = "'0.03137"
$ws.Range("D13").Value = "'0.02906"
$ws.Range("D14").Value = "'0.09385"
$ws.Range("D15").Value = "'3.757"
$ws.Range("D16").Value = "'0.001525"
$ws.Range("D18").Value = "'0.0005962"
$ws.Range("D19").Value = "'0.006081"
$ws.Range("D20").Value = "'0.001235"

$ws.Range("D21").Value = "'0.004613"
$ws.Range("E21").Value = "20HotbitTokenHTBBestin24h"

$ws.Range("D22").Value = "'0.00006902"
$ws.Range("D23").Value = "'3.500"
$ws.Range("D24").Value = "'2.150"
$ws.Range("D25").Value = "'0.3195"

$ws.Range("D40").Value = "'0.03644"

$ws.Range("D41").Value = "'0.006272"
$ws.Range("E41").Value = "40KickTokenKICK"

$ws.Range("D42").Value = "'0.1053"
$ws.Range("D43").Value = "'0.002771"
$ws.Range("D44").Value = "'0.008491"
$ws.Range("D45").Value = "'0.00005281"
$ws.Range("D48").Value = "'0.002286"
$ws.Range("D49").Value = "'0.00002101"
$ws.Range("D50").Value = "'0.0002001"
